$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2741.4285
$ws.Range("I40").Value = 1297.5
$ws.Range("J40").Value = 4666.6665
$ws.Range("K40").Value = 1297.5
$ws.Range("L40").Value = 4666.6665
$ws.Range("M40").Value = -1122.5
$ws.Range("N40").Value = -5016.6665
$ws.Range("H43").Value = 5970.75
$ws.Range("I43").Value = 5236.375
$ws.Range("K43").Value = 5236.375
$ws.Range("M43").Value = -5167.375
$ws.Range("H53").Value = 945.75
$ws.Range("I53").Value = 152.4
$ws.Range("J53").Value = 1306.3636
$ws.Range("K53").Value = 152.4
$ws.Range("L53").Value = 1306.3636
$ws.Range("M53").Value = 484.6
$ws.Range("N53").Value = -2580.3636
$ws.Range("H70").Value = 4680274.5
$ws.Range("J70").Value = 57525.145
$ws.Range("L70").Value = 172575.435
$ws.Range("N70").Value = -173115.435
$ws.Range("H73").Value = 4680274.5
$ws.Range("J73").Value = 57525.145
$ws.Range("L73").Value = 172575.435
$ws.Range("N73").Value = -174447.435
$ws.Range("H74").Value = 4733.0835
$ws.Range("I74").Value = 4527
$ws.Range("K74").Value = 4527
$ws.Range("M74").Value = -3591
$ws.Range("H77").Value = 4733.0835
$ws.Range("I77").Value = 4527
$ws.Range("K77").Value = 22635
$ws.Range("M77").Value = -17955
$ws.Range("H80").Value = 3109314
$ws.Range("I80").Value = 7255490
$ws.Range("J80").Value = 345196.66
$ws.Range("K80").Value = 21766470
$ws.Range("L80").Value = 1035589.98
$ws.Range("M80").Value = -21765472
$ws.Range("N80").Value = -1037585.98
$ws.Range("H83").Value = 3109314
$ws.Range("I83").Value = 7255490
$ws.Range("J83").Value = 345196.66
$ws.Range("K83").Value = 65299410
$ws.Range("L83").Value = 3106769.94
$ws.Range("M83").Value = -65294418
$ws.Range("N83").Value = -3116753.94
$ws.Range("H106").Value = 4535.364
$ws.Range("I106").Value = 4488.3
$ws.Range("K106").Value = 4488.3
$ws.Range("M106").Value = -3857.3
$ws.Range("H107").Value = 1939
$ws.Range("J107").Value = 2717.8333
$ws.Range("L107").Value = 2717.8333
$ws.Range("N107").Value = -6557.8333
$ws.Range("H113").Value = 4014.1428
$ws.Range("J113").Value = 4820
$ws.Range("L113").Value = 4820
$ws.Range("N113").Value = -11328
$ws.Range("H137").Value = 8649.965
$ws.Range("I137").Value = 12408.517
$ws.Range("K137").Value = 37225.551
$ws.Range("M137").Value = -34675.551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4281.12
$ws.Range("I32").Value = 4281.12
$ws.Range("K32").Value = 4281.12
$ws.Range("M32").Value = -3994.12
$ws.Range("H63").Value = 4550.5625
$ws.Range("I63").Value = 4101.25
$ws.Range("K63").Value = 4101.25
$ws.Range("M63").Value = -3415.25
$ws.Range("H66").Value = 4550.5625
$ws.Range("I66").Value = 4101.25
$ws.Range("K66").Value = 20506.25
$ws.Range("M66").Value = -17074.25
$ws.Range("H102").Value = 3035.4285
$ws.Range("I102").Value = 2569.5264
$ws.Range("J102").Value = 4019
$ws.Range("K102").Value = 2569.5264
$ws.Range("L102").Value = 4019
$ws.Range("M102").Value = -947.5264000000002
$ws.Range("N102").Value = -7263

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H139").Value = 129997.664
$ws.Range("J139").Value = 129997.664
$ws.Range("L139").Value = 129997.664
$ws.Range("N139").Value = -140277.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1482.6666
$ws.Range("I94").Value = 1734.25
$ws.Range("J94").Value = 1281.4
$ws.Range("K94").Value = 1734.25
$ws.Range("L94").Value = 1281.4
$ws.Range("M94").Value = -1283.25
$ws.Range("N94").Value = -2183.4
$ws.Range("H99").Value = 4759.6665
$ws.Range("I99").Value = 4617.364
$ws.Range("J99").Value = 4916.2
$ws.Range("K99").Value = 4617.364
$ws.Range("L99").Value = 4916.2
$ws.Range("M99").Value = -3119.364
$ws.Range("N99").Value = -7912.2
$ws.Range("H122").Value = 12865.857
$ws.Range("I122").Value = 22570.5
$ws.Range("J122").Value = 4043.4546
$ws.Range("K122").Value = 67711.5
$ws.Range("L122").Value = 12130.3638
$ws.Range("M122").Value = -65261.5
$ws.Range("N122").Value = -17030.3638
$ws.Range("H126").Value = 4759.6665
$ws.Range("I126").Value = 4617.364
$ws.Range("J126").Value = 4916.2
$ws.Range("K126").Value = 13852.092
$ws.Range("L126").Value = 14748.6
$ws.Range("M126").Value = -11382.092
$ws.Range("N126").Value = -19688.6
$ws.Range("H132").Value = 15022.837
$ws.Range("I132").Value = 15368.619
$ws.Range("K132").Value = 46105.857
$ws.Range("M132").Value = -43575.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 346.5
$ws.Range("J2").Value = 138.83333
$ws.Range("L2").Value = 832.9999799999999
$ws.Range("N2").Value = -1058.99998
$ws.Range("H107").Value = 981.1905
$ws.Range("I107").Value = 1272.9231
$ws.Range("J107").Value = 507.125
$ws.Range("K107").Value = 3818.7693
$ws.Range("L107").Value = 1521.375
$ws.Range("M107").Value = -1898.7693
$ws.Range("N107").Value = -5361.375
$ws.Range("H113").Value = 2203.0645
$ws.Range("I113").Value = 2101.75
$ws.Range("J113").Value = 2267.0527
$ws.Range("K113").Value = 6305.25
$ws.Range("L113").Value = 6801.158100000001
$ws.Range("M113").Value = -4135.25
$ws.Range("N113").Value = -11141.1581

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 164.64285
$ws.Range("I2").Value = 224.66667
$ws.Range("J2").Value = 56.6
$ws.Range("K2").Value = 224.66667
$ws.Range("L2").Value = 56.6
$ws.Range("M2").Value = -111.66667
$ws.Range("N2").Value = -282.6
$ws.Range("H20").Value = 20323.6
$ws.Range("J20").Value = 20323.6
$ws.Range("L20").Value = 20323.6
$ws.Range("N20").Value = -20813.6
$ws.Range("H96").Value = 60261
$ws.Range("J96").Value = 60261
$ws.Range("L96").Value = 60261
$ws.Range("N96").Value = -65753
$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("N109").Value = -29080
$ws.Range("H113").Value = 2429.8125
$ws.Range("I113").Value = 2489.4285
$ws.Range("J113").Value = 2012.5
$ws.Range("K113").Value = 2489.4285
$ws.Range("L113").Value = 2012.5
$ws.Range("M113").Value = -319.4285
$ws.Range("N113").Value = -6352.5
$ws.Range("H122").Value = 3584.1667
$ws.Range("I122").Value = 3579.4546
$ws.Range("J122").Value = 3597.125
$ws.Range("K122").Value = 10738.3638
$ws.Range("L122").Value = 10791.375
$ws.Range("M122").Value = -8288.363799999999
$ws.Range("N122").Value = -15691.375
$ws.Range("H126").Value = 2565.611
$ws.Range("I126").Value = 2536.5
$ws.Range("J126").Value = 2798.5
$ws.Range("K126").Value = 7609.5
$ws.Range("L126").Value = 8395.5
$ws.Range("M126").Value = -5139.5
$ws.Range("N126").Value = -13335.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1197.9048
$ws.Range("I16").Value = 1177.6428
$ws.Range("K16").Value = 1177.6428
$ws.Range("M16").Value = -1007.6428
$ws.Range("H68").Value = 6145.8335
$ws.Range("I68").Value = 5218.75
$ws.Range("J68").Value = 8000
$ws.Range("K68").Value = 5218.75
$ws.Range("L68").Value = 8000
$ws.Range("M68").Value = -4469.75
$ws.Range("N68").Value = -9498
$ws.Range("H71").Value = 6145.8335
$ws.Range("I71").Value = 5218.75
$ws.Range("J71").Value = 8000
$ws.Range("K71").Value = 26093.75
$ws.Range("L71").Value = 40000
$ws.Range("M71").Value = -22349.75
$ws.Range("N71").Value = -47488
$ws.Range("H122").Value = 14787.77
$ws.Range("I122").Value = 14787.77
$ws.Range("K122").Value = 44363.31
$ws.Range("M122").Value = -41913.31

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8231.916999999999
$ws.Range("I62").Value = 9966.333000000001
$ws.Range("J62").Value = 7653.778
$ws.Range("K62").Value = 9966.333000000001
$ws.Range("L62").Value = 7653.778
$ws.Range("M62").Value = -9342.333000000001
$ws.Range("N62").Value = -8901.778
$ws.Range("H65").Value = 8231.916999999999
$ws.Range("I65").Value = 9966.333000000001
$ws.Range("J65").Value = 7653.778
$ws.Range("K65").Value = 49831.665
$ws.Range("L65").Value = 38268.89
$ws.Range("M65").Value = -46711.665
$ws.Range("N65").Value = -44508.89
$ws.Range("H96").Value = 1766.6666
$ws.Range("I96").Value = 1473.2142
$ws.Range("J96").Value = 2793.75
$ws.Range("K96").Value = 1473.2142
$ws.Range("L96").Value = 2793.75
$ws.Range("M96").Value = -100.2141999999999
$ws.Range("N96").Value = -5539.75
$ws.Range("H122").Value = 11937798
$ws.Range("I122").Value = 15192478
$ws.Range("J122").Value = 3970.5557
$ws.Range("K122").Value = 45577434
$ws.Range("L122").Value = 11911.6671
$ws.Range("M122").Value = -45574984
$ws.Range("N122").Value = -16811.6671
